# Timpi Sortare.xlsx - "Add files via upload" edit
# Adds a Radix-Sort-base-10 column (F/G) of timings to Sheet1, renames/fixes
# the radix sort row, adds two brand-new result tables on Sheet2 and Sheet3,
# and leaves Sheet2 as the active tab/sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# Sheet1 ("Sheet1") - extra Radix-Sort-baza-10 column + row14 fix
# ---------------------------------------------------------------------------

# New timings in columns F/G for the existing test rows
$ws1.Range("F7").Value  = 0.166
$ws1.Range("G7").Value  = 0.092

$ws1.Range("F8").Value  = 0.067
$ws1.Range("G8").Value  = 0.046

$ws1.Range("F9").Value  = 0.057
$ws1.Range("G9").Value  = 0.032

$ws1.Range("F10").Value = 0.046
$ws1.Range("G10").Value = 0.041

$ws1.Range("F13").Value = 0.012
$ws1.Range("G13").Value = 0.01

$ws1.Range("F15").Value = 0.062
$ws1.Range("G15").Value = 0.035

# Row 14 used to be the broken radixSort row (all "    -" placeholders and a
# "crapa" note); it now has real numbers plus a renamed label.
$ws1.Range("B14").Value = "Radix Sort, baza 10"
$ws1.Range("C14").Value = 0.05
$ws1.Range("D14").Value = 0.073
$ws1.Range("E14").Value = "CRAPA"
$ws1.Range("F14").Value = 0.088
$ws1.Range("G14").Value = 0.045
$ws1.Range("M14").Value = "(radixSort e in pluton, pe testele mici)"

# The old standalone warning in M15 no longer applies - drop it.
$ws1.Range("M15").ClearContents()

# ---------------------------------------------------------------------------
# Sheet3 - new "Teste aproape sortate 200k-500k" results table
# (populated before Sheet2's own new table so that new shared strings are
# appended in the same first-seen order as the original edit.)
# ---------------------------------------------------------------------------

$ws3.Columns.Item(3).ColumnWidth = 30.71

$ws3.Range("C3:F3").HorizontalAlignment = -4108
$ws3.Range("C3").Value = 'Teste "aproape" sortate 200k-500k'
$ws3.Range("C3:F3").Merge()

$ws3.Range("D5").Value = "Heapsort"
$ws3.Range("E5").Value = "Test 1"
$ws3.Range("F5").Value = "Test 2"
$ws3.Range("G5").Value = "Test 3"
$ws3.Range("H5").Value = "Test 4"

$ws3.Range("C6").Value = "Heapsort"
$ws3.Range("D6").Value = 0.3129
$ws3.Range("E6").Value = 0.4842
$ws3.Range("F6").Value = 0.7342

$ws3.Range("C7").Value = "Mergesort"
$ws3.Range("D7").Value = 0.1249
$ws3.Range("E7").Value = 0.1093
$ws3.Range("F7").Value = 0.1874

$ws3.Range("C8").Value = "Quicksort(pivot la final)"
$ws3.Range("D8").Value = 0.4998
$ws3.Range("E8").Value = 2.7237
$ws3.Range("F8").Value = 2.3232

$ws3.Range("C9").Value = "Quicksort(pivot la mijloc)"
$ws3.Range("D9").Value = 0.0468
$ws3.Range("E9").Value = 0.0937
$ws3.Range("F9").Value = 0.1405

$ws3.Range("C10").Value = "Bubblesort"
$ws3.Range("D10").Value = "    - "
$ws3.Range("E10").Value = "    - "
$ws3.Range("F10").Value = "    - "
$ws3.Range("G10").Value = "    - "
$ws3.Range("H10").Value = "    - "

$ws3.Range("C11").Value = "Shellsort"
$ws3.Range("D11").Value = "    -"
$ws3.Range("E11").Value = "    -"
$ws3.Range("F11").Value = "    -"
$ws3.Range("G11").Value = "    -"
$ws3.Range("H11").Value = "    -"

$ws3.Range("C12").Value = "Counting Sort"
$ws3.Range("D12").Value = 0.0156
$ws3.Range("E12").Value = 0.0156
$ws3.Range("F12").Value = 0.0156

$ws3.Range("C13").Value = "Radix Sort"
$ws3.Range("D13").Value = 0.0781
$ws3.Range("E13").Value = 0.1562
$ws3.Range("F13").Value = 0.203

$ws3.Range("C14").Value = "std::sort"
$ws3.Range("D14").Value = 0.0624
$ws3.Range("E14").Value = 0.0624
$ws3.Range("F14").Value = 0.0937

# ---------------------------------------------------------------------------
# Sheet2 - new "Teste Random, 1-3M" results table
# ---------------------------------------------------------------------------

$ws2.Range("C3:F3").HorizontalAlignment = -4108
$ws2.Range("C3:F3").Merge()

$ws2.Range("E5:G5").HorizontalAlignment = -4108
$ws2.Range("E5").Font.Name = "Calibri"
$ws2.Range("E5").Value = "Teste Random, 1-3M"
$ws2.Range("E5:G5").Merge()

$ws2.Range("D7").Value = "Heapsort"
$ws2.Range("E7").Value = "Test 1"
$ws2.Range("F7").Value = "Test 2"
$ws2.Range("G7").Value = "Test 3"

$ws2.Range("D8").Value = "Heapsort"
$ws2.Range("E8").Value = 2.0307
$ws2.Range("F8").Value = 2.8118
$ws2.Range("G8").Value = 4.3114

$ws2.Range("D9").Value = "Mergesort"
$ws2.Range("E9").Value = 0.8904
$ws2.Range("F9").Value = 1.031
$ws2.Range("G9").Value = 1.4371

$ws2.Range("D10").Value = "Quicksort(pivot la final)"
$ws2.Range("E10").Value = 0.5155
$ws2.Range("F10").Value = 0.7186
$ws2.Range("G10").Value = 1.2028

$ws2.Range("D11").Value = "Quicksort(pivot la mijloc)"
$ws2.Range("E11").Value = 0.5936
$ws2.Range("F11").Value = 0.7498
$ws2.Range("G11").Value = 1.2497

$ws2.Range("D12").Value = "Bubblesort"
$ws2.Range("E12").Value = "    - "
$ws2.Range("F12").Value = "    - "
$ws2.Range("G12").Value = "    - "

$ws2.Range("D13").Value = "Shellsort"
$ws2.Range("E13").Value = "    -"
$ws2.Range("F13").Value = "    -"
$ws2.Range("G13").Value = "    -"

$ws2.Range("D14").Value = "Counting Sort"
$ws2.Range("E14").Value = 0.0312
$ws2.Range("F14").Value = 0.0468
$ws2.Range("G14").Value = 0.0624

$ws2.Range("D15").Value = "Radix Sort"
$ws2.Range("E15").Value = 0.5467
$ws2.Range("F15").Value = 0.6561
$ws2.Range("G15").Value = 0.9841

$ws2.Range("D16").Value = "std::sort"
$ws2.Range("E16").Value = 0.5155
$ws2.Range("F16").Value = 0.7654
$ws2.Range("G16").Value = 1.1716

$ws2.Range("D6:G6").Font.Name = "Calibri"
$ws2.Range("D3:G16").Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# Selections / active sheet (Sheet2 becomes the active tab)
# ---------------------------------------------------------------------------

$ws1.Range("H12").Select()
$ws3.Range("G7").Select()
$ws2.Activate()
$ws2.Range("H8").Select()
